$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two brochure links to point at the new "region-2" folder
# instead of the old "sismic" folder (same file names/paths otherwise).
$ws.Range("B2").Value = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/offshore/region-2/aibt/AIBT_Courses_Fees_2021_VOL_2.2.pdf"
$ws.Range("B3").Value = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/offshore/region-2/aibt/AIBTSISMIC_Q4_Brochure_1OCT-31DEC21_VOL1.1.pdf"

# Turn those two cells into real clickable hyperlinks (this is also what
# introduces the "Hyperlink" cell style / underline font used below).
$ws.Hyperlinks.Add($ws.Range("B2"), $ws.Range("B2").Value2)
$ws.Hyperlinks.Add($ws.Range("B3"), $ws.Range("B3").Value2)

# Widen column B so the longer URLs still fit.
$ws.Columns.Item(2).ColumnWidth = 160

# Move the active selection to B4.
$ws.Range("B4").Select()
